$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Alex Lockwood, Wisconsin, 1998-10-21, Manager, Checking
$ws.Range("B2").Value = "Alex"
$ws.Range("C2").Value = "Lockwood"
$ws.Range("D2").Value = "Wisconsin"
$ws.Range("E2").Value = Get-Date -Year 1998 -Month 10 -Day 21 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("F2").Value = "Manager"
$ws.Range("G2").Value = "Checking"

# Update row 3: Tyler Johnson, Chicago, 1990-12-12, Analyst, Savings
$ws.Range("B3").Value = "Tyler"
$ws.Range("C3").Value = "Johnson"
$ws.Range("D3").Value = "Chicago"
$ws.Range("E3").Value = Get-Date -Year 1990 -Month 12 -Day 12 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("F3").Value = "Analyst"
$ws.Range("G3").Value = "Savings"

# Update the selected cell in the sheet view
$ws.Range("F8").Select()
